$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Range("B16").Value = 6221766
$ws.Range("F16").Value = "Kairat Almaty"
$ws.Range("G16").Value = "FK Kaspyi Aktau"
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = "H"
$ws.Range("K16").Value = 1.55
$ws.Range("L16").Value = 3.8
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = 1.65
$ws.Range("O16").Value = 4
$ws.Range("P16").Value = 4.5
$ws.Range("Q16").Value = -0.75
$ws.Range("R16").Value = 1.8
$ws.Range("S16").Value = 2
$ws.Range("T16").Value = 2.75
$ws.Range("U16").Value = 1.925
$ws.Range("V16").Value = 1.875
$ws.Range("W16").Value = 0.6499999999999999
$ws.Range("X16").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = 0.8
$ws.Range("AA16").Value = -1
$ws.Range("AB16").Value = 0.925
$ws.Range("AC16").Value = -1

# Row 17
$ws.Range("B17").Value = 6221640
$ws.Range("F17").Value = "Ordabasy"
$ws.Range("G17").Value = "Zhetysu"
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "H"
$ws.Range("K17").Value = 1.333
$ws.Range("L17").Value = 4.333
$ws.Range("M17").Value = 8
$ws.Range("N17").Value = 1.333
$ws.Range("O17").Value = 4.333
$ws.Range("P17").Value = 8
$ws.Range("Q17").Value = -1.25
$ws.Range("R17").Value = 1.75
$ws.Range("S17").Value = 1.95
$ws.Range("T17").Value = 2.5
$ws.Range("U17").Value = 1.9
$ws.Range("V17").Value = 1.9
$ws.Range("W17").Value = 0.333
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = -1
$ws.Range("Z17").Value = 0.75
$ws.Range("AA17").Value = -1
$ws.Range("AB17").Value = -1
$ws.Range("AC17").Value = 0.8999999999999999

# Row 18
$ws.Range("B18").Value = 6221641
$ws.Range("F18").Value = "Tobol Kostanay"
$ws.Range("G18").Value = "Shakhter Karagandy"
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = "H"
$ws.Range("K18").Value = 1.4
$ws.Range("L18").Value = 4.333
$ws.Range("M18").Value = 6
$ws.Range("N18").Value = 1.333
$ws.Range("O18").Value = 4.75
$ws.Range("P18").Value = 6.5
$ws.Range("Q18").Value = -1.5
$ws.Range("R18").Value = 1.925
$ws.Range("S18").Value = 1.875
$ws.Range("T18").Value = 2.75
$ws.Range("U18").Value = 1.75
$ws.Range("V18").Value = 1.95
$ws.Range("W18").Value = 0.333
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Z18").Value = -1
$ws.Range("AA18").Value = 0.875
$ws.Range("AB18").Value = 0.375
$ws.Range("AC18").Value = -0.5

# Row 19
$ws.Range("B19").Value = 6221642
$ws.Range("F19").Value = "Kaisar Kyzylorda"
$ws.Range("G19").Value = "FK Aksu"
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = "D"
$ws.Range("K19").Value = 2.45
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = 2.7
$ws.Range("N19").Value = 2.15
$ws.Range("O19").Value = 3.3
$ws.Range("P19").Value = 3.1
$ws.Range("Q19").Value = -0.25
$ws.Range("R19").Value = 1.9
$ws.Range("S19").Value = 1.9
$ws.Range("T19").Value = 2.5
$ws.Range("U19").Value = 1.925
$ws.Range("V19").Value = 1.875
$ws.Range("W19").Value = -1
$ws.Range("X19").Value = 2.3
$ws.Range("Y19").Value = -1
$ws.Range("Z19").Value = -0.5
$ws.Range("AA19").Value = 0.45
$ws.Range("AB19").Value = 0.925
$ws.Range("AC19").Value = -1

# Row 20
$ws.Range("B20").Value = 6221639
$ws.Range("F20").Value = "FC Astana"
$ws.Range("G20").Value = "FK Aktobe"
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = "A"
$ws.Range("K20").Value = 1.45
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 5.75
$ws.Range("N20").Value = 1.533
$ws.Range("O20").Value = 3.75
$ws.Range("P20").Value = 5
$ws.Range("Q20").Value = -1
$ws.Range("R20").Value = 1.9
$ws.Range("S20").Value = 1.9
$ws.Range("T20").Value = 2.75
$ws.Range("U20").Value = 1.95
$ws.Range("V20").Value = 1.85
$ws.Range("W20").Value = -1
$ws.Range("X20").Value = -1
$ws.Range("Y20").Value = 4
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 0.8999999999999999
$ws.Range("AB20").Value = 0.95
$ws.Range("AC20").Value = -1

# Row 31
$ws.Range("B31").Value = 6221771
$ws.Range("F31").Value = "FK Maktaaral"
$ws.Range("G31").Value = "FK Kaspyi Aktau"
$ws.Range("H31").Value = 3
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = "H"
$ws.Range("K31").Value = 2.05
$ws.Range("L31").Value = 3.3
$ws.Range("M31").Value = 3.1
$ws.Range("N31").Value = 2.1
$ws.Range("O31").Value = 3.25
$ws.Range("P31").Value = 3.1
$ws.Range("Q31").Value = -0.25
$ws.Range("R31").Value = 1.875
$ws.Range("S31").Value = 1.925
$ws.Range("T31").Value = 2.25
$ws.Range("U31").Value = 1.875
$ws.Range("V31").Value = 1.925
$ws.Range("W31").Value = 1.1
$ws.Range("X31").Value = -1
$ws.Range("Y31").Value = -1
$ws.Range("Z31").Value = 0.875
$ws.Range("AA31").Value = -1
$ws.Range("AB31").Value = 0.875
$ws.Range("AC31").Value = -1

# Row 32
$ws.Range("B32").Value = 6221772
$ws.Range("F32").Value = "FK Kyzylzhar"
$ws.Range("G32").Value = "FK Atyrau"
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = "H"
$ws.Range("K32").Value = 1.909
$ws.Range("L32").Value = 3.3
$ws.Range("M32").Value = 3.5
$ws.Range("N32").Value = 1.909
$ws.Range("O32").Value = 3.4
$ws.Range("P32").Value = 3.4
$ws.Range("Q32").Value = -0.5
$ws.Range("R32").Value = 1.975
$ws.Range("S32").Value = 1.825
$ws.Range("T32").Value = 2.25
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 1.8
$ws.Range("W32").Value = 0.909
$ws.Range("X32").Value = -1
$ws.Range("Y32").Value = -1
$ws.Range("Z32").Value = 0.9750000000000001
$ws.Range("AA32").Value = -1
$ws.Range("AB32").Value = -1
$ws.Range("AC32").Value = 0.8

# Row 66
$ws.Range("B66").Value = 6221674
$ws.Range("F66").Value = "Zhetysu"
$ws.Range("G66").Value = "FK Atyrau"
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 3
$ws.Range("J66").Value = "A"
$ws.Range("K66").Value = 2.5
$ws.Range("L66").Value = 3.2
$ws.Range("M66").Value = 2.5
$ws.Range("N66").Value = 2.5
$ws.Range("O66").Value = 3.2
$ws.Range("P66").Value = 2.5
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = 1.9
$ws.Range("S66").Value = 1.9
$ws.Range("T66").Value = 2.25
$ws.Range("U66").Value = 1.925
$ws.Range("V66").Value = 1.875
$ws.Range("W66").Value = -1
$ws.Range("X66").Value = -1
$ws.Range("Y66").Value = 1.5
$ws.Range("Z66").Value = -1
$ws.Range("AA66").Value = 0.8999999999999999
$ws.Range("AB66").Value = 0.925
$ws.Range("AC66").Value = -1

# Row 67
$ws.Range("B67").Value = 6221673
$ws.Range("F67").Value = "Shakhter Karagandy"
$ws.Range("G67").Value = "FK Aksu"
$ws.Range("H67").Value = 2
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = "H"
$ws.Range("K67").Value = 2.1
$ws.Range("L67").Value = 3.25
$ws.Range("M67").Value = 3
$ws.Range("N67").Value = 2.1
$ws.Range("O67").Value = 3.3
$ws.Range("P67").Value = 3
$ws.Range("Q67").Value = -0.25
$ws.Range("R67").Value = 1.925
$ws.Range("S67").Value = 1.875
$ws.Range("T67").Value = 2.5
$ws.Range("U67").Value = 1.975
$ws.Range("V67").Value = 1.825
$ws.Range("W67").Value = 1.1
$ws.Range("X67").Value = -1
$ws.Range("Y67").Value = -1
$ws.Range("Z67").Value = 0.925
$ws.Range("AA67").Value = -1
$ws.Range("AB67").Value = 0.9750000000000001
$ws.Range("AC67").Value = -1

# Row 95
$ws.Range("B95").Value = 6221693
$ws.Range("F95").Value = "Zhetysu"
$ws.Range("G95").Value = "Shakhter Karagandy"
$ws.Range("H95").Value = 1
$ws.Range("I95").Value = 3
$ws.Range("J95").Value = "A"
$ws.Range("K95").Value = 2
$ws.Range("L95").Value = 3.4
$ws.Range("M95").Value = 3.1
$ws.Range("N95").Value = 2.2
$ws.Range("O95").Value = 3.3
$ws.Range("P95").Value = 2.8
$ws.Range("Q95").Value = -0.25
$ws.Range("R95").Value = 1.95
$ws.Range("S95").Value = 1.85
$ws.Range("T95").Value = 2.5
$ws.Range("U95").Value = 1.85
$ws.Range("V95").Value = 1.95
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = 1.8
$ws.Range("Z95").Value = -1
$ws.Range("AA95").Value = 0.8500000000000001
$ws.Range("AB95").Value = 0.8500000000000001
$ws.Range("AC95").Value = -1

# Row 96
$ws.Range("B96").Value = 6221698
$ws.Range("F96").Value = "FK Maktaaral"
$ws.Range("G96").Value = "FK Aktobe"
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = "A"
$ws.Range("K96").Value = 4.333
$ws.Range("L96").Value = 3.5
$ws.Range("M96").Value = 1.666
$ws.Range("N96").Value = 4.2
$ws.Range("O96").Value = 3.4
$ws.Range("P96").Value = 1.7
$ws.Range("Q96").Value = 0.75
$ws.Range("R96").Value = 1.825
$ws.Range("S96").Value = 1.975
$ws.Range("T96").Value = 2.5
$ws.Range("U96").Value = 1.925
$ws.Range("V96").Value = 1.875
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = 0.7
$ws.Range("Z96").Value = -0.5
$ws.Range("AA96").Value = 0.4875
$ws.Range("AB96").Value = 0.925
$ws.Range("AC96").Value = -1

# Row 104
$ws.Range("B104").Value = 6221703
$ws.Range("F104").Value = "Shakhter Karagandy"
$ws.Range("G104").Value = "FK Aktobe"
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = "A"
$ws.Range("K104").Value = 3.6
$ws.Range("L104").Value = 3.5
$ws.Range("M104").Value = 1.8
$ws.Range("N104").Value = 3.1
$ws.Range("O104").Value = 3.5
$ws.Range("P104").Value = 1.909
$ws.Range("Q104").Value = 0.5
$ws.Range("R104").Value = 1.825
$ws.Range("S104").Value = 1.975
$ws.Range("T104").Value = 2.5
$ws.Range("U104").Value = 1.75
$ws.Range("V104").Value = 1.95
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = 0.909
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 0.9750000000000001
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.95

# Row 105
$ws.Range("B105").Value = 6221699
$ws.Range("F105").Value = "FK Maktaaral"
$ws.Range("G105").Value = "Kaisar Kyzylorda"
$ws.Range("H105").Value = 2
$ws.Range("I105").Value = 2
$ws.Range("J105").Value = "D"
$ws.Range("K105").Value = 3.1
$ws.Range("L105").Value = 3.2
$ws.Range("M105").Value = 2.1
$ws.Range("N105").Value = 2.1
$ws.Range("O105").Value = 3.1
$ws.Range("P105").Value = 3.2
$ws.Range("Q105").Value = -0.25
$ws.Range("R105").Value = 1.85
$ws.Range("S105").Value = 1.95
$ws.Range("T105").Value = 2.25
$ws.Range("U105").Value = 1.975
$ws.Range("V105").Value = 1.725
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = 2.1
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = -0.5
$ws.Range("AA105").Value = 0.475
$ws.Range("AB105").Value = 0.9750000000000001
$ws.Range("AC105").Value = -1

# Row 177
$ws.Range("B177").Value = 6221753
$ws.Range("F177").Value = "FK Aksu"
$ws.Range("G177").Value = "Tobol Kostanay"
$ws.Range("H177").Value = 0
$ws.Range("I177").Value = 3
$ws.Range("J177").Value = "A"
$ws.Range("K177").Value = 2.75
$ws.Range("L177").Value = 3.1
$ws.Range("M177").Value = 2.375
$ws.Range("N177").Value = 2.625
$ws.Range("O177").Value = 3.2
$ws.Range("P177").Value = 2.45
$ws.Range("Q177").Value = 0
$ws.Range("R177").Value = 2
$ws.Range("S177").Value = 1.8
$ws.Range("T177").Value = 2.5
$ws.Range("U177").Value = 1.9
$ws.Range("V177").Value = 1.9
$ws.Range("W177").Value = -1
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = 1.45
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.8
$ws.Range("AB177").Value = 0.8999999999999999
$ws.Range("AC177").Value = -1

# Row 178
$ws.Range("B178").Value = 6221815
$ws.Range("F178").Value = "FK Atyrau"
$ws.Range("G178").Value = "Kairat Almaty"
$ws.Range("H178").Value = 0
$ws.Range("I178").Value = 0
$ws.Range("J178").Value = "D"
$ws.Range("K178").Value = 3
$ws.Range("L178").Value = 3
$ws.Range("M178").Value = 2.25
$ws.Range("N178").Value = 3.1
$ws.Range("O178").Value = 3.1
$ws.Range("P178").Value = 2.15
$ws.Range("Q178").Value = 0.25
$ws.Range("R178").Value = 1.85
$ws.Range("S178").Value = 1.95
$ws.Range("T178").Value = 2.25
$ws.Range("U178").Value = 1.8
$ws.Range("V178").Value = 2
$ws.Range("W178").Value = -1
$ws.Range("X178").Value = 2.1
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = 0.425
$ws.Range("AA178").Value = -0.5
$ws.Range("AB178").Value = -1
$ws.Range("AC178").Value = 1

# Row 179
$ws.Range("B179").Value = 6221752
$ws.Range("F179").Value = "FK Kyzylzhar"
$ws.Range("G179").Value = "Kaisar Kyzylorda"
$ws.Range("H179").Value = 0
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = "A"
$ws.Range("K179").Value = 1.833
$ws.Range("L179").Value = 3.2
$ws.Range("M179").Value = 4
$ws.Range("N179").Value = 1.85
$ws.Range("O179").Value = 3.2
$ws.Range("P179").Value = 4
$ws.Range("Q179").Value = -0.5
$ws.Range("R179").Value = 1.9
$ws.Range("S179").Value = 1.9
$ws.Range("T179").Value = 2
$ws.Range("U179").Value = 1.775
$ws.Range("V179").Value = 2.025
$ws.Range("W179").Value = -1
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = 3
$ws.Range("Z179").Value = -1
$ws.Range("AA179").Value = 0.8999999999999999
$ws.Range("AB179").Value = -1
$ws.Range("AC179").Value = 1.025

# Row 180
$ws.Range("B180").Value = 6221814
$ws.Range("F180").Value = "Okzhetpes Kokshetau"
$ws.Range("G180").Value = "FK Maktaaral"
$ws.Range("H180").Value = 1
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = "D"
$ws.Range("K180").Value = 2.3
$ws.Range("L180").Value = 3.1
$ws.Range("M180").Value = 2.8
$ws.Range("N180").Value = 2.3
$ws.Range("O180").Value = 3.1
$ws.Range("P180").Value = 2.8
$ws.Range("Q180").Value = 0
$ws.Range("R180").Value = 1.75
$ws.Range("S180").Value = 2.05
$ws.Range("T180").Value = 2.25
$ws.Range("U180").Value = 1.875
$ws.Range("V180").Value = 1.925
$ws.Range("W180").Value = -1
$ws.Range("X180").Value = 2.1
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0
$ws.Range("AA180").Value = -0
$ws.Range("AB180").Value = -0.5
$ws.Range("AC180").Value = 0.4625

# Row 181
$ws.Range("B181").Value = 6221755
$ws.Range("F181").Value = "FK Aktobe"
$ws.Range("G181").Value = "Zhetysu"
$ws.Range("H181").Value = 2
$ws.Range("I181").Value = 2
$ws.Range("J181").Value = "D"
$ws.Range("K181").Value = 1.727
$ws.Range("L181").Value = 3.5
$ws.Range("M181").Value = 4
$ws.Range("N181").Value = 1.444
$ws.Range("O181").Value = 4
$ws.Range("P181").Value = 6
$ws.Range("Q181").Value = -1.25
$ws.Range("R181").Value = 1.975
$ws.Range("S181").Value = 1.825
$ws.Range("T181").Value = 2.75
$ws.Range("U181").Value = 1.8
$ws.Range("V181").Value = 2
$ws.Range("W181").Value = -1
$ws.Range("X181").Value = 3
$ws.Range("Y181").Value = -1
$ws.Range("Z181").Value = -1
$ws.Range("AA181").Value = 0.825
$ws.Range("AB181").Value = 0.8
$ws.Range("AC181").Value = -1

# Row 182
$ws.Range("B182").Value = 6221816
$ws.Range("F182").Value = "FK Kaspyi Aktau"
$ws.Range("G182").Value = "Ordabasy"
$ws.Range("H182").Value = 1
$ws.Range("I182").Value = 0
$ws.Range("J182").Value = "H"
$ws.Range("K182").Value = 3.4
$ws.Range("L182").Value = 3.4
$ws.Range("M182").Value = 1.909
$ws.Range("N182").Value = 4.2
$ws.Range("O182").Value = 4
$ws.Range("P182").Value = 1.571
$ws.Range("Q182").Value = 0.75
$ws.Range("R182").Value = 1.95
$ws.Range("S182").Value = 1.75
$ws.Range("T182").Value = 3
$ws.Range("U182").Value = 1.975
$ws.Range("V182").Value = 1.825
$ws.Range("W182").Value = 3.2
$ws.Range("X182").Value = -1
$ws.Range("Y182").Value = -1
$ws.Range("Z182").Value = 0.95
$ws.Range("AA182").Value = -1
$ws.Range("AB182").Value = -1
$ws.Range("AC182").Value = 0.825
